$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:J1 (reuse the existing header style from E1) ---
$ws.Range("F1:J1").Value = "x"
$ws.Range("E1").Copy()
$ws.Range("F1:J1").PasteSpecial(-4122)

$ws.Range("F1").Value = "FreundesListe"
$ws.Range("G1").Value = "Fitness Level (X/100)"
$ws.Range("H1").Value = "Läuft los"
$ws.Range("I1").Value = "Lauf Start Position"
$ws.Range("J1").Value = "Zeit des Laufs"

# --- Row 2: new data for the existing user (Mezix) ---
$ws.Range("F2").Value = "Phine;Mezix2;"
$ws.Range("F2").HorizontalAlignment = 1

$ws.Range("G2").Value = 50

$ws.Range("H2").Formula = "=TRUE()"
$ws.Range("H2").NumberFormat = '"TRUE";"TRUE";"FALSE"'

# --- New rows: friends added as their own user rows ---
$ws.Range("A3").Value = "Mezix2"
$ws.Range("B3").Value = "password"
$ws.Range("C3").Value = "Felix"
$ws.Range("D3").Value = "Swimmer"
$ws.Range("E3").Value = "21"

$ws.Range("A4").Value = "Phine"

$ws.Range("A5").Value = "Enric"

# --- Column widths ---
$ws.Columns.Item(6).ColumnWidth = 17.08666666666667
$ws.Columns.Item(7).ColumnWidth = 22.51666666666667
$ws.Columns.Item(8).ColumnWidth = 14.846666666666666
$ws.Columns.Item(9).ColumnWidth = 21.526666666666667
$ws.Columns.Item(10).ColumnWidth = 20.726666666666667

# --- Selection ---
$ws.Range("F2").Select()
